$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old header row (X1) and the trailing extra row (chose):
# shift "Guitar"/"Valli"/"Blubbeladf" up so they occupy A1:A3.
$ws.Range("A1").Value = "Guitar"
$ws.Range("A2").Value = "Valli"
$ws.Range("A3").Value = "Blubbeladf"

# Drop the now-unused rows 4 and 5 entirely so the sheet only spans A1:A3.
$ws.Rows.Item(4).Resize(2).Delete()

# Apply a Text number format to the remaining cells.
$ws.Range("A1:A3").NumberFormat = "@"

# Reselect the A1:A3 block (mirrors the author ending up with A1:A3 selected).
$ws.Range("A1:A3").Select()
